# "10Th - MB for single stock and added new group"
# MarketBeat rank tracker: a new weekly snapshot ("Jun_26"/"Jun_27") is
# inserted in front of the existing date columns, and a new ratings group
# ("Benchmark" / "Evercore ISI") is appended as new rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three new columns in front of column B. This shifts the
#    existing B:E date columns (and any per-cell styling, such as the
#    "latest change" highlight fills) three columns to the right, landing
#    on E:H - exactly mirroring how Excel natively keeps data/format
#    anchored to the same logical column when columns are inserted.
$ws.Columns("B:D").Insert()

# 2) New header dates for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# 3) Default the new snapshot columns to "UN" (unchanged) for every
#    existing analyst row (rows 2-27).
$ws.Range("B2:D27").Value = "UN"

# Goldman Sachs Group (row 6) got a new rating on 6/26/2018 that shows up
# across the three new columns.
$ws.Range("B6:D6").Value = "6/26/2018,Reiterates,Sell,`$195.00"

# 4) Brand new analyst rows for the newly tracked group.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
